$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Coin/Link text updates (row reshuffle for Maker / EnergySwap / Aptos)
$textUpdates = @{
    'B47' = 'Maker'
    'C47' = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
    'B48' = 'EnergySwap'
    'C48' = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
    'B49' = 'Aptos'
    'C49' = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
}
foreach ($addr in $textUpdates.Keys) {
    $ws.Range($addr).Value = $textUpdates[$addr]
}

# Price (column D) updates - force Text storage so values like "1.000" /
# "0.4680" keep their trailing zeros instead of being read as numbers.
$priceUpdates = @{
    'D2' = '30.349.53'
    'D3' = '1.868.29'
    'D4' = '1.000'
    'D5' = '236.28'
    'D7' = '0.4680'
    'D8' = '0.2849'
    'D9' = '0.06545'
    'D10' = '21.93'
    'D11' = '0.07930'
    'D12' = '97.76'
    'D13' = '1.874.79'
    'D14' = '5.164'
    'D15' = '0.6814'
    'D16' = '281.14'
    'D17' = '30.343.71'
    'D18' = '13.25'
    'D19' = '1.000'
    'D20' = '5.426'
    'D21' = '2.116.99'
    'D22' = '0.000007331'
    'D23' = '1.000'
    'D24' = '6.165'
    'D25' = '166.24'
    'D27' = '19.14'
    'D28' = '1.945'
    'D29' = '1.389'
    'D30' = '0.09786'
    'D31' = '4.416'
    'D32' = '1.483'
    'D33' = '4.115'
    'D34' = '0.04742'
    'D36' = '0.7115'
    'D37' = '2.717'
    'D38' = '0.01873'
    'D39' = '6.372'
    'D40' = '2.558'
    'D41' = '74.93'
    'D43' = '0.8530'
    'D44' = '0.4212'
    'D45' = '1.0000'
    'D46' = '103.64'
    'D47' = '971.07'
    'D48' = '9.421'
    'D49' = '7.242'
    'D50' = '34.23'
    'D51' = '0.1131'
}
foreach ($addr in $priceUpdates.Keys) {
    $ws.Range($addr).NumberFormat = '@'
    $ws.Range($addr).Value = $priceUpdates[$addr]
    $ws.Range($addr).ClearFormats()
}

# Volume(1h) (column E) updates
$volumeUpdates = @{
    'E2' = '  +0.82%  '
    'E3' = '  +0.68%  '
    'E4' = '  +0.03%  '
    'E5' = '  +0.95%  '
    'E6' = '  -0.01%  '
    'E7' = '  -0.17%  '
    'E8' = '  +1.08%  '
    'E9' = '  -0.11%  '
    'E10' = '  +8.43%  '
    'E11' = '  +1.73%  '
    'E12' = '  +0.57%  '
    'E13' = '  +1.00%  '
    'E14' = '  +1.75%  '
    'E15' = '  +1.78%  '
    'E16' = '  -0.80%  '
    'E17' = '  +0.73%  '
    'E18' = '  +5.25%  '
    'E19' = '  +0.00%  '
    'E20' = '  +0.47%  '
    'E21' = '  +1.04%  '
    'E22' = '  +1.38%  '
    'E23' = '  +0.07%  '
    'E24' = '  +0.50%  '
    'E25' = '  -0.88%  '
    'E26' = '  -1.12%  '
    'E27' = '  +0.77%  '
    'E28' = '  +0.92%  '
    'E29' = '  +3.98%  '
    'E30' = '  +1.28%  '
    'E31' = '  +0.71%  '
    'E32' = '  +0.93%  '
    'E33' = '  +0.44%  '
    'E34' = '  +1.59%  '
    'E35' = '  +4.99%  '
    'E36' = '  +1.97%  '
    'E37' = '  +0.61%  '
    'E38' = '  +0.82%  '
    'E39' = '  +1.10%  '
    'E40' = '  +2.19%  '
    'E41' = '  +4.73%  '
    'E42' = '  +2.47%  '
    'E43' = '  -0.85%  '
    'E44' = '  +1.18%  '
    'E45' = '  +0.01%  '
    'E46' = '  -0.60%  '
    'E47' = '  -4.94%  '
    'E48' = '  +2.12%  '
    'E49' = '  -0.13%  '
    'E50' = '  +1.35%  '
}
foreach ($addr in $volumeUpdates.Keys) {
    $ws.Range($addr).Value = $volumeUpdates[$addr]
}
